# Update countries & provincias Spain
# Applies the COVID-19 data refresh scraped at 19 de Junio de 2020, 06:45
# (replacing the previous 05:28 snapshot): some countries' totals were
# updated and, as a consequence of re-ranking, a handful of adjacent rows
# swapped country labels while keeping their row position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $values) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($row, $i + 1).Value = $values[$i]
    }
}

# --- Timestamp banner (row 1) ---
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 19 de Junio de 2020 a las 06:45"

# --- Row 56: Kazajistan (label unchanged, totals refreshed) ---
Set-Row 56 @("Kazajistan", 16351, 474, 10139, 6107, 0, 5, 105)

# --- Rows 80-81: Haiti overtakes Guinea ---
Set-Row 80 @("Haiti", 4916, 228, 24, 4808, 0, 2, 84)
Set-Row 81 @("Guinea", 4841, 0, 3467, 1348, 0, 0, 26)

# --- Rows 95-96: Kirguistan overtakes Somalia ---
Set-Row 95 @("Kirguistan", 2789, 132, 1961, 796, 0, 1, 32)
Set-Row 96 @("Somalia", 2719, 0, 724, 1907, 0, 0, 88)

# --- Row 164: Mongolia (label unchanged, totals refreshed) ---
Set-Row 164 @("Mongolia", 204, 3, 132, 72, 0, 0, 0)

# --- Rows 202-203: Fiyi overtakes Dominica (values tied, only labels swap) ---
$ws.Cells.Item(202, 1).Value = "Fiyi"
$ws.Cells.Item(203, 1).Value = "Dominica"

# --- Rows 208-211: Santa Sede / Islas Turcas y Caicos / Montserrat / Seychelles reshuffle ---
Set-Row 208 @("Santa Sede", 12, 0, 12, 0, 0, 0, 0)
Set-Row 209 @("Islas Turcas y Caicos", 12, 0, 11, 0, 0, 0, 1)
Set-Row 210 @("Montserrat", 11, 0, 10, 0, 0, 0, 1)
Set-Row 211 @("Seychelles", 11, 0, 11, 0, 0, 0, 0)
